$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9, shifting existing rows 9..98 down to 10..99
$ws.Rows.Item(9).Insert()

# Populate the new row 9 with this week's data
$ws.Cells.Item(9, 1).Value = 3
$ws.Cells.Item(9, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(9, 3).Value = "Coquimbo"
$ws.Cells.Item(9, 4).Value = 44545
$ws.Cells.Item(9, 5).Value = 5
$ws.Cells.Item(9, 6).Value = 100112052
$ws.Cells.Item(9, 7).Value = "Albahaca"
$ws.Cells.Item(9, 8).Value = "Sin especificar"
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 105
$ws.Cells.Item(9, 11).Value = 4000
$ws.Cells.Item(9, 12).Value = 4500
$ws.Cells.Item(9, 13).Value = 4238
$ws.Cells.Item(9, 14).Value = "`$/docena de matas"
$ws.Cells.Item(9, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(9, 16).Value = 706
$ws.Cells.Item(9, 17).Value = 6
$ws.Cells.Item(9, 18).Value = "Hortaliza"
